$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 70, shifting existing rows 70-113 down to 71-114.
$ws.Rows.Item(70).Insert()

# Populate the newly inserted row 70 with the new weekly record.
$ws.Cells.Item(70, 1).Value = 2
$ws.Cells.Item(70, 2).Value = 'Comercializadora del Agro de Limarí'
$ws.Cells.Item(70, 3).Value = 'Coquimbo'
$ws.Cells.Item(70, 4).Value = 44609
$ws.Cells.Item(70, 5).Value = 4
$ws.Cells.Item(70, 6).Value = 100112024
$ws.Cells.Item(70, 7).Value = 'Choclo'
$ws.Cells.Item(70, 8).Value = 'Choclero'
$ws.Cells.Item(70, 9).Value = 'Primera'
$ws.Cells.Item(70, 10).Value = 50000
$ws.Cells.Item(70, 11).Value = 130
$ws.Cells.Item(70, 12).Value = 150
$ws.Cells.Item(70, 13).Value = 140
$ws.Cells.Item(70, 14).Value = '$/unidad'
$ws.Cells.Item(70, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(70, 16).Value = 140
$ws.Cells.Item(70, 17).Value = 1
$ws.Cells.Item(70, 18).Value = 'Hortaliza'
